$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3's formatting (styles + cell types) into row 4, then
# overwrite with the actual trade values.
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

$ws.Cells.Item(4, 1).Value = 42641.54247685185
$ws.Cells.Item(4, 2).Value = $true
$ws.Cells.Item(4, 3).Value = 9948.4599999999991
$ws.Cells.Item(4, 4).Value = 9942
$ws.Cells.Item(4, 5).Value = 79.319999999999993
$ws.Cells.Item(4, 6).Value = 79.22
$ws.Cells.Item(4, 7).Value = $true
$ws.Cells.Item(4, 8).Value = -0.13
$ws.Cells.Item(4, 9).Value = $false
